$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-08 Sunday" "2024-09-09 Monday"

Replace-Text "837×6=" "821×5="
Replace-Text "332×8=" "393×8="
Replace-Text "454×4=" "228×8="
Replace-Text "620×6=" "916×5="
Replace-Text "341×5=" "450×7="
Replace-Text "686×9=" "934×4="
Replace-Text "740×5=" "431×3="
Replace-Text "376×3=" "385×3="
Replace-Text "111×5=" "212×9="
Replace-Text "762×4=" "150×6="
Replace-Text "899×8=" "994×5="
Replace-Text "630×2=" "662×2="
Replace-Text "963×3=" "839×3="
Replace-Text "391×7=" "386×3="
Replace-Text "434×5=" "510×9="
Replace-Text "828×2=" "761×2="
Replace-Text "464×7=" "837×8="
Replace-Text "344×2=" "115×5="
Replace-Text "582×7=" "773×5="
Replace-Text "214×6=" "697×6="
Replace-Text "853×7=" "225×9="
Replace-Text "607×6=" "501×3="
Replace-Text "294×3=" "589×5="
Replace-Text "900×5=" "935×2="
Replace-Text "226×6=" "102×2="
